$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh -- updates Price (D) and Volume(1h) (E) columns,
# plus a 3-row reorder/update at the bottom (FTXToken/ARBITRUM/MXToken).

$ws.Range("D2").Value = "37.006.78"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "1.978.05"
$ws.Range("E3").Value = "  -3.23%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.89"
$ws.Range("E5").Value = "  -6.04%  "

$ws.Range("E6").Value = "  -4.30%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.06"
$ws.Range("E8").Value = "  -6.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.00"
$ws.Range("E9").Value = "  +3.20%  "

$ws.Range("E10").Value = "  -5.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0746"
$ws.Range("E11").Value = "  -7.28%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0984"
$ws.Range("E12").Value = "  -4.90%  "

$ws.Range("D13").Value = "2.264.80"
$ws.Range("E13").Value = "  -3.41%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "13.90"
$ws.Range("E14").Value = "  -6.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.79"
$ws.Range("E15").Value = "  -3.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.748"
$ws.Range("E16").Value = "  -9.61%  "

$ws.Range("E17").Value = "  -7.68%  "

$ws.Range("D18").Value = "1.964.94"
$ws.Range("E18").Value = "  -3.99%  "

$ws.Range("D19").Value = "36.918.98"
$ws.Range("E19").Value = "  -1.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.94"
$ws.Range("E20").Value = "  -3.31%  "

$ws.Range("D21").Value = "0.0₃0804"
$ws.Range("E21").Value = "  -6.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.77"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -6.78%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  -11.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.55"
$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.55"
$ws.Range("E28").Value = "  -7.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.93"
$ws.Range("E29").Value = "  -5.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.125"
$ws.Range("E30").Value = "  -11.20%  "

$ws.Range("E31").Value = "  -6.21%  "

$ws.Range("E32").Value = "  -3.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.37"
$ws.Range("E33").Value = "  -9.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0610"
$ws.Range("E34").Value = "  -8.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("E35").Value = "  -6.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("E36").Value = "  -7.92%  "

$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.77"
$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("E39").Value = "  -6.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.10"
$ws.Range("E40").Value = "  -5.86%  "

$ws.Range("E41").Value = "  -0.50%  "

$ws.Range("D42").Value = "1.411.68"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("E43").Value = "  -7.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0889"
$ws.Range("E44").Value = "  -8.66%  "

$ws.Range("E45").Value = "  -7.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "86.79"
$ws.Range("E46").Value = "  -5.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.21"
$ws.Range("E47").Value = "  -7.62%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.86"
$ws.Range("E48").Value = "  +18.53%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.992"
$ws.Range("E49").Value = "  -6.21%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.59"
$ws.Range("E51").Value = "  -11.86%  "
